# Natmi following Dr Hou advice
# Rebuild the Col4a3 -> Cd47 sending/target cluster cross-table (ECs/FAPs/sCs x ECs/FAPs/sCs)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Col4a3", "Cd47", "ECs", 2, 0.6666666666666666, 0.07529266666666666, 0.225878, 0.008636301470765068, 0.008636301470765068, 3, 1, 59.45197733333333, 178.355932, 0.304222453049858, 0.304222453049858, 4.476297912032888, 40.286681208296, 0.002627356818714245, 0.002627356818714245),
    @("ECs", "Col4a3", "Cd47", "FAPs", 2, 0.6666666666666666, 0.07529266666666666, 0.225878, 0.008636301470765068, 0.008636301470765068, 3, 1, 107.1770123333333, 321.531037, 0.548436823552382, 0.5484368235523819, 8.069643063942888, 72.62678757548599, 0.004736465745867158, 0.004736465745867158),
    @("ECs", "Col4a3", "Cd47", "sCs", 2, 0.6666666666666666, 0.07529266666666666, 0.225878, 0.008636301470765068, 0.008636301470765068, 3, 1, 28.793724, 86.38117199999999, 0.1473407233977601, 0.1473407233977601, 2.167956263224, 19.511606369016, 0.001272478906183664, 0.001272478906183664),
    @("FAPs", "Col4a3", "Cd47", "ECs", 3, 1, 7.575405333333333, 22.726216, 0.868922394680866, 0.868922394680866, 3, 1, 59.45197733333333, 178.355932, 0.304222453049858, 0.304222453049858, 450.3728261681458, 4053.355435513312, 0.2643457024197699, 0.2643457024197699),
    @("FAPs", "Col4a3", "Cd47", "FAPs", 3, 1, 7.575405333333333, 22.726216, 0.868922394680866, 0.868922394680866, 3, 1, 107.1770123333333, 321.531037, 0.548436823552382, 0.5484368235523819, 811.9093108406657, 7307.183797565992, 0.4765490380523033, 0.4765490380523032),
    @("FAPs", "Col4a3", "Cd47", "sCs", 3, 1, 7.575405333333333, 22.726216, 0.868922394680866, 0.868922394680866, 3, 1, 28.793724, 86.38117199999999, 0.1473407233977601, 0.1473407233977601, 218.124130356128, 1963.117173205152, 0.1280276542087928, 0.1280276542087928),
    @("sCs", "Col4a3", "Cd47", "ECs", 3, 1, 1.067463, 3.202389, 0.1224413038483689, 0.1224413038483689, 3, 1, 59.45197733333333, 178.355932, 0.304222453049858, 0.304222453049858, 63.462786080172, 571.165074721548, 0.03724939381137381, 0.03724939381137381),
    @("sCs", "Col4a3", "Cd47", "FAPs", 3, 1, 1.067463, 3.202389, 0.1224413038483689, 0.1224413038483689, 3, 1, 107.1770123333333, 321.531037, 0.548436823552382, 0.5484368235523819, 114.407495116377, 1029.667456047393, 0.06715131975421151, 0.06715131975421149),
    @("sCs", "Col4a3", "Cd47", "sCs", 3, 1, 1.067463, 3.202389, 0.1224413038483689, 0.1224413038483689, 3, 1, 28.793724, 86.38117199999999, 0.1473407233977601, 0.1473407233977601, 30.736235002212, 276.626115019908, 0.01804059028278363, 0.01804059028278362)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value2 = $val
        $c++
    }
    $r++
}
